# Refresh cryptos list values (price/volume columns D & E, and the
# Chainlink/Uniswap row swap) to match the latest scrape, per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Numeric-looking text in column D (e.g. "218.43") is written with a
# leading apostrophe so Excel keeps storing it as text (matching the
# source data, which treats prices as plain strings) instead of
# auto-converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.806.93"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.641.41"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.18%  "

# Row 5: BNB
$ws.Range("D5").Value = "'218.43"
$ws.Range("E5").Value = "  +0.68%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = "  -0.21%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.14%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  -0.16%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  -0.92%  "

# Row 10: Solana
$ws.Range("E10").Value = "  -0.08%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +0.41%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.870.07"

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.640.44"
$ws.Range("E13").Value = "  -0.15%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  -0.65%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  -0.56%  "

# Row 16: Litecoin
$ws.Range("E16").Value = "  +0.86%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.821.48"
$ws.Range("E17").Value = "  +0.13%  "

# Row 18: ShibaInu
$ws.Range("E18").Value = "  -0.84%  "

# Row 19: BitcoinCash
$ws.Range("D19").Value = "'216.29"

# Row 20: Dai
$ws.Range("E20").Value = "  -0.17%  "

# Row 21: Uniswap (was Chainlink)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'4.36"
$ws.Range("E21").Value = "  -0.10%  "

# Row 22: Chainlink (was Uniswap)
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'6.60"
$ws.Range("E22").Value = "  +4.91%  "

# Row 23: Toncoin
$ws.Range("E23").Value = "  -1.97%  "

# Row 24: Avalanche
$ws.Range("D24").Value = "'9.17"
$ws.Range("E24").Value = "  -2.23%  "

# Row 25: Monero
$ws.Range("D25").Value = "'147.47"
$ws.Range("E25").Value = "  +1.67%  "

# Row 26: BinanceUSD
$ws.Range("E26").Value = "  -0.23%  "

# Row 27: Stellar
$ws.Range("D27").Value = "'0.119"
$ws.Range("E27").Value = "  +0.04%  "

# Row 28: Cosmos
$ws.Range("E28").Value = "  +0.16%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'15.74"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30: Hedera
$ws.Range("D30").Value = "'0.0510"
$ws.Range("E30").Value = "  -0.84%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  +1.01%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  +1.77%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.89%  "

# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  +0.63%  "

# Row 35: Maker
$ws.Range("D35").Value = "1.264.37"
$ws.Range("E35").Value = "  -2.22%  "

# Row 36: HuobiToken
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  +0.31%  "

# Row 37: VeChain
$ws.Range("E37").Value = "  -0.10%  "

# Row 38: ImmutableX
$ws.Range("E38").Value = "  -2.08%  "

# Row 39: ARBITRUM
$ws.Range("D39").Value = "'0.816"
$ws.Range("E39").Value = "  -1.46%  "

# Row 40: PaxDollar
$ws.Range("E40").Value = "  -0.15%  "

# Row 41: TrustWalletToken
$ws.Range("D41").Value = "'0.804"
$ws.Range("E41").Value = "  -0.82%  "

# Row 42: FraxShare
$ws.Range("E42").Value = "  -0.69%  "

# Row 43: RocketPoolETH
$ws.Range("D43").Value = "1.780.37"
$ws.Range("E43").Value = "  -0.72%  "

# Row 44: MXToken
$ws.Range("E44").Value = "  -4.40%  "

# Row 45: Quant
$ws.Range("E45").Value = "  +1.02%  "

# Row 46: Aave
$ws.Range("D46").Value = "'60.90"
$ws.Range("E46").Value = "  +0.68%  "

# Row 47: RenderToken
$ws.Range("D47").Value = "'1.60"
$ws.Range("E47").Value = "  -0.39%  "

# Row 48: BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  +7.34%  "

# Row 49: Cronos
$ws.Range("D49").Value = "'0.0515"
$ws.Range("E49").Value = "  -0.69%  "

# Row 50: EnergySwap
$ws.Range("D50").Value = "'7.56"
$ws.Range("E50").Value = "  -1.58%  "

# Row 51: Algorand
$ws.Range("E51").Value = "  -1.68%  "
